# Apply updated crypto price/volume figures as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.975.94"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").Value = "2.405.16"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'562.94"
$ws.Range("E5").Value = "  +1.65%  "

$ws.Range("D6").Value = "'142.44"
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = "  -0.56%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("E10").Value = "  -1.80%  "

$ws.Range("D11").Value = "'5.28"
$ws.Range("E11").Value = "  -2.13%  "

$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").Value = "'25.59"
$ws.Range("E13").Value = "  -2.55%  "

$ws.Range("D14").Value = "'0.0000173"
$ws.Range("E14").Value = "  -0.86%  "

$ws.Range("D15").Value = "2.839.55"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("D16").Value = "61.895.62"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").Value = "2.404.13"
$ws.Range("E17").Value = "  -0.48%  "

$ws.Range("D18").Value = "'11.24"
$ws.Range("E18").Value = "  +1.26%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'322.30"
$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.84"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("E21").Value = "  -1.11%  "

$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").Value = "'65.94"
$ws.Range("E23").Value = "  +1.54%  "

$ws.Range("D24").Value = "'1.76"
$ws.Range("E24").Value = "  -0.68%  "

$ws.Range("D25").Value = "'8.81"
$ws.Range("E25").Value = "  -4.51%  "

$ws.Range("D26").Value = "'568.03"
$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.518.46"
$ws.Range("E28").Value = "  -0.90%  "

$ws.Range("D29").Value = "0.0₃0940"
$ws.Range("E29").Value = "  +0.73%  "

$ws.Range("D30").Value = "'8.20"
$ws.Range("E30").Value = "  -1.91%  "

$ws.Range("E31").Value = "  -2.81%  "

$ws.Range("E32").Value = "  -0.46%  "

$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("E34").Value = "  -2.45%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").Value = "'4.68"
$ws.Range("E36").Value = "  -2.76%  "

$ws.Range("D37").Value = "'5.45"
$ws.Range("E37").Value = "  -5.04%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.380"
$ws.Range("E38").Value = "  -1.09%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'151.07"
$ws.Range("E39").Value = "  +2.72%  "

$ws.Range("E40").Value = "  -1.00%  "

$ws.Range("D41").Value = "'1.80"
$ws.Range("E41").Value = "  -8.75%  "

$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("D43").Value = "'2.26"
$ws.Range("E43").Value = "  -1.45%  "

$ws.Range("D44").Value = "'147.83"
$ws.Range("E44").Value = "  -2.87%  "

$ws.Range("E45").Value = "  -0.56%  "

$ws.Range("E46").Value = "  -2.77%  "

$ws.Range("D47").Value = "'19.87"
$ws.Range("E47").Value = "  -2.78%  "

$ws.Range("D48").Value = "'0.588"
$ws.Range("E48").Value = "  -0.34%  "

$ws.Range("D49").Value = "'0.0916"
$ws.Range("E49").Value = "  +0.55%  "

$ws.Range("E50").Value = "  -0.96%  "

$ws.Range("D51").Value = "'11.53"
$ws.Range("E51").Value = "  +0.56%  "
